$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("27-10-2021", 2.88, 3.79, 4.68, 5.82, 0.75),
    @("28-10-2021", 2.81, 3.84, 4.92, 5.82, 0.82),
    @("29-10-2021", 2.8, 3.59, 4.92, 5.82, 0.63)
)

$startRow = 206
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
}
